$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo: "Current SectionMapping" -> "Current Section Mapping"
$ws.Range("C1").Value = "Current Section Mapping"

# Reset the view: scroll back to column A (remove topLeftCell="D1") and
# move the active selection to C1 (single cell).
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C1").Select()
